$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 69.84614999999999
$ws.Range("I9").Value = 81.55556
$ws.Range("J9").Value = 43.5
$ws.Range("K9").Value = 81.55556
$ws.Range("L9").Value = 43.5
$ws.Range("M9").Value = 87.44444
$ws.Range("N9").Value = -381.5

$ws.Range("H19").Value = 354.33334
$ws.Range("I19").Value = 374.9091
$ws.Range("J19").Value = 331.7
$ws.Range("K19").Value = 374.9091
$ws.Range("L19").Value = 331.7
$ws.Range("M19").Value = -199.9091
$ws.Range("N19").Value = -681.7

$ws.Range("H70").Value = 3461.25
$ws.Range("J70").Value = 3969
$ws.Range("L70").Value = 11907
$ws.Range("N70").Value = -12447

$ws.Range("H73").Value = 3461.25
$ws.Range("J73").Value = 3969
$ws.Range("L73").Value = 11907
$ws.Range("N73").Value = -13779

$ws.Range("H76").Value = 8000
$ws.Range("I76").Value = 8000
$ws.Range("K76").Value = 8000
$ws.Range("M76").Value = -7685

$ws.Range("H79").Value = 8000
$ws.Range("I79").Value = 8000
$ws.Range("K79").Value = 8000
$ws.Range("M79").Value = -6908

$ws.Range("H87").Value = 58462.332
$ws.Range("J87").Value = 67754.8
$ws.Range("L87").Value = 67754.8
$ws.Range("N87").Value = -70250.8

$ws.Range("H90").Value = 58462.332
$ws.Range("J90").Value = 67754.8
$ws.Range("L90").Value = 203264.4
$ws.Range("N90").Value = -215744.4

$ws.Range("H94").Value = 7939.385
$ws.Range("I94").Value = 7939.385
$ws.Range("K94").Value = 7939.385
$ws.Range("M94").Value = -7488.385

$ws.Range("H125").Value = 2444.3333
$ws.Range("I125").Value = 2444.3333
$ws.Range("K125").Value = 21998.9997
$ws.Range("M125").Value = -19538.9997

$ws.Range("H132").Value = 13845.777
$ws.Range("I132").Value = 13701.5625
$ws.Range("K132").Value = 41104.6875
$ws.Range("M132").Value = -38574.6875

$ws.Range("H138").Value = 2072.875
$ws.Range("I138").Value = 1916.8
$ws.Range("J138").Value = 2333
$ws.Range("K138").Value = 5750.4
$ws.Range("L138").Value = 6999
$ws.Range("M138").Value = -610.3999999999996
$ws.Range("N138").Value = -17279

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 2703.8
$ws.Range("J46").Value = 2625.1667
$ws.Range("L46").Value = 2625.1667
$ws.Range("N46").Value = -3263.1667

$ws.Range("H110").Value = 2908.08
$ws.Range("I110").Value = 1205.8334
$ws.Range("K110").Value = 1205.8334
$ws.Range("M110").Value = 839.1666

$ws.Range("H122").Value = 1938.8334
$ws.Range("I122").Value = 1938.8334
$ws.Range("K122").Value = 5816.5002
$ws.Range("M122").Value = -3366.5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 219.52632
$ws.Range("I80").Value = 173.16667
$ws.Range("J80").Value = 240.92308
$ws.Range("K80").Value = 173.16667
$ws.Range("L80").Value = 240.92308
$ws.Range("M80").Value = 824.8333299999999
$ws.Range("N80").Value = -2236.92308

$ws.Range("H83").Value = 219.52632
$ws.Range("I83").Value = 173.16667
$ws.Range("J83").Value = 240.92308
$ws.Range("K83").Value = 865.8333500000001
$ws.Range("L83").Value = 1204.6154
$ws.Range("M83").Value = 4126.16665
$ws.Range("N83").Value = -11188.6154

$ws.Range("H99").Value = 1738.2
$ws.Range("I99").Value = 1996.625
$ws.Range("K99").Value = 1996.625
$ws.Range("M99").Value = -498.625

$ws.Range("H105").Value = 1713.5
$ws.Range("I105").Value = 2030
$ws.Range("J105").Value = 1397
$ws.Range("K105").Value = 2030
$ws.Range("L105").Value = 1397
$ws.Range("M105").Value = -283
$ws.Range("N105").Value = -4891

$ws.Range("H132").Value = 150000
$ws.Range("J132").Value = 150000
$ws.Range("L132").Value = 150000
$ws.Range("N132").Value = -160120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1670.3334
$ws.Range("I16").Value = 1554.125
$ws.Range("J16").Value = 2600
$ws.Range("K16").Value = 1554.125
$ws.Range("L16").Value = 2600
$ws.Range("M16").Value = -1267.125
$ws.Range("N16").Value = -3174

$ws.Range("H62").Value = 3200
$ws.Range("J62").Value = 3000
$ws.Range("L62").Value = 3000
$ws.Range("N62").Value = -4248

$ws.Range("H65").Value = 3200
$ws.Range("J65").Value = 3000
$ws.Range("L65").Value = 15000
$ws.Range("N65").Value = -21240

$ws.Range("H113").Value = 1670.3334
$ws.Range("I113").Value = 1554.125
$ws.Range("J113").Value = 2600
$ws.Range("K113").Value = 1554.125
$ws.Range("L113").Value = 2600
$ws.Range("M113").Value = 615.875
$ws.Range("N113").Value = -6940

$ws.Range("H132").Value = 4821.6665
$ws.Range("I132").Value = 4128.5713
$ws.Range("K132").Value = 12385.7139
$ws.Range("M132").Value = -9855.713899999999

$ws.Range("H134").Value = 4335
$ws.Range("I134").Value = 1199
$ws.Range("K134").Value = 3597
$ws.Range("M134").Value = -1062

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 169267.67
$ws.Range("I109").Value = 202010.2
$ws.Range("K109").Value = 606030.6000000001
$ws.Range("M109").Value = -604990.6000000001

$ws.Range("H131").Value = 1750
$ws.Range("J131").Value = 2000
$ws.Range("L131").Value = 6000
$ws.Range("N131").Value = -16080

$ws.Range("H137").Value = 2466.6667
$ws.Range("I137").Value = 2000
$ws.Range("J137").Value = 2700
$ws.Range("K137").Value = 6000
$ws.Range("L137").Value = 8100
$ws.Range("M137").Value = -900
$ws.Range("N137").Value = -18300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8074.125
$ws.Range("I7").Value = 6998.3335
$ws.Range("K7").Value = 6998.3335
$ws.Range("M7").Value = -6886.3335

$ws.Range("H9").Value = 130.85715
$ws.Range("I9").Value = 130.85715
$ws.Range("K9").Value = 130.85715
$ws.Range("M9").Value = 93.14285000000001

$ws.Range("H16").Value = 699.3333
$ws.Range("I16").Value = 359
$ws.Range("J16").Value = 1124.75
$ws.Range("K16").Value = 359
$ws.Range("L16").Value = 1124.75
$ws.Range("M16").Value = -189
$ws.Range("N16").Value = -1464.75

$ws.Range("H61").Value = 3496.7144
$ws.Range("I61").Value = 2437.1177
$ws.Range("K61").Value = 2437.1177
$ws.Range("M61").Value = -2235.1177

$ws.Range("H113").Value = 3496.7144
$ws.Range("I113").Value = 2437.1177
$ws.Range("K113").Value = 2437.1177
$ws.Range("M113").Value = -267.1176999999998

$ws.Range("H126").Value = 8074.125
$ws.Range("I126").Value = 6998.3335
$ws.Range("K126").Value = 20995.0005
$ws.Range("M126").Value = -18525.0005

$ws.Range("H132").Value = 4361.926
$ws.Range("I132").Value = 3332
$ws.Range("K132").Value = 9996
$ws.Range("M132").Value = -7466

$ws.Range("H136").Value = 6799.8
$ws.Range("I136").Value = 5000
$ws.Range("K136").Value = 15000
$ws.Range("M136").Value = -12450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3983.3333
$ws.Range("I136").Value = 2138.5
$ws.Range("J136").Value = 6197.1333
$ws.Range("K136").Value = 6415.5
$ws.Range("L136").Value = 18591.3999
$ws.Range("M136").Value = -3865.5
$ws.Range("N136").Value = -23691.3999
